$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.035182332627914
$ws.Range("D2").Value = 1.029357660681784
$ws.Range("E2").Value = 1.04324360465346
$ws.Range("F2").Value = 1.051913823851382
$ws.Range("I2").Value = 1.03134161344135
$ws.Range("J2").Value = 1.040297505370909
$ws.Range("K2").Value = 1.0321718053833
$ws.Range("L2").Value = 1.04601801888951
$ws.Range("M2").Value = 1.054664007687564
$ws.Range("N2").Value = 1.017106412322584

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036942325711251
$ws.Range("D3").Value = 1.029750866589934
$ws.Range("E3").Value = 1.044799191421415
$ws.Range("F3").Value = 1.053534451870518
$ws.Range("I3").Value = 1.031410473938938
$ws.Range("J3").Value = 1.041696850780066
$ws.Range("K3").Value = 1.032374257636609
$ws.Range("L3").Value = 1.047382654709337
$ws.Range("M3").Value = 1.056095302517026
$ws.Range("N3").Value = 1.017594728209205

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038078744999993
$ws.Range("D4").Value = 1.030004897780531
$ws.Range("E4").Value = 1.04580363544453
$ws.Range("F4").Value = 1.054580688442747
$ws.Range("I4").Value = 1.031453268851647
$ws.Range("J4").Value = 1.042599666925593
$ws.Range("K4").Value = 1.032504118621964
$ws.Range("L4").Value = 1.048263057974017
$ws.Range("M4").Value = 1.057018555124507
$ws.Range("N4").Value = 1.017909171142915

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038555932402221
$ws.Range("D5").Value = 1.030111595620222
$ws.Range("E5").Value = 1.046225406724308
$ws.Range("F5").Value = 1.055019959229197
$ws.Range("I5").Value = 1.031470837888524
$ws.Range("J5").Value = 1.042978586886911
$ws.Range("K5").Value = 1.032558438855329
$ws.Range("L5").Value = 1.048632566056704
$ws.Range("M5").Value = 1.057406009307275
$ws.Range("N5").Value = 1.018041000429604

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038636021659413
$ws.Range("D6").Value = 1.030129504938905
$ws.Range("E6").Value = 1.046296195171217
$ws.Range("F6").Value = 1.055093681815884
$ws.Range("I6").Value = 1.03147376306346
$ws.Range("J6").Value = 1.043042172956733
$ws.Range("K6").Value = 1.032567543415666
$ws.Range("L6").Value = 1.048694572427501
$ws.Range("M6").Value = 1.057471024936542
$ws.Range("N6").Value = 1.018063113997706

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038085123394628
$ws.Range("D7").Value = 1.030006323864434
$ws.Range("E7").Value = 1.045809273109797
$ws.Range("F7").Value = 1.054586560212371
$ws.Range("I7").Value = 1.031453505267842
$ws.Range("J7").Value = 1.042604732508675
$ws.Range("K7").Value = 1.032504845526229
$ws.Range("L7").Value = 1.048267997754737
$ws.Range("M7").Value = 1.057023734967981
$ws.Range("N7").Value = 1.017910934071522

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035777639190834
$ws.Range("D8").Value = 1.029490628382574
$ws.Range("E8").Value = 1.04376976944836
$ws.Range("F8").Value = 1.052462030520538
$ws.Range("I8").Value = 1.031365250087661
$ws.Range("J8").Value = 1.040770977104939
$ws.Range("K8").Value = 1.03224046058015
$ws.Range("L8").Value = 1.046479750327819
$ws.Range("M8").Value = 1.055148325613929
$ws.Range("N8").Value = 1.017271760069284

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03169241042364
$ws.Range("D9").Value = 1.028578903965162
$ws.Range("E9").Value = 1.040159105992702
$ws.Range("F9").Value = 1.048699284259999
$ws.Range("I9").Value = 1.031196234843219
$ws.Range("J9").Value = 1.037518826077261
$ws.Range("K9").Value = 1.031765881887206
$ws.Range("L9").Value = 1.04330817565782
$ws.Range("M9").Value = 1.051820984343072
$ws.Range("N9").Value = 1.016133564133688

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028955126326753
$ws.Range("D10").Value = 1.027969141535362
$ws.Range("E10").Value = 1.03773996518119
$ws.Range("F10").Value = 1.046177226014938
$ws.Range("I10").Value = 1.031074480267779
$ws.Range("J10").Value = 1.035335991767315
$ws.Range("K10").Value = 1.031443680964606
$ws.Range("L10").Value = 1.041179358562246
$ws.Range("M10").Value = 1.049586852339097
$ws.Range("N10").Value = 1.015366531946505

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027766372494171
$ws.Range("D11").Value = 1.027704660988068
$ws.Range("E11").Value = 1.036689431373512
$ws.Range("F11").Value = 1.045081762372516
$ws.Range("I11").Value = 1.031019604969636
$ws.Range("J11").Value = 1.034387144089246
$ws.Range("K11").Value = 1.031302791644956
$ws.Range("L11").Value = 1.040253980972282
$ws.Range("M11").Value = 1.048615521645865
$ws.Range("N11").Value = 1.015032391157962

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027324273891058
$ws.Range("D12").Value = 1.027606354783216
$ws.Range("E12").Value = 1.036298746798029
$ws.Range("F12").Value = 1.044674333623874
$ws.Range("I12").Value = 1.030998897942262
$ws.Range("J12").Value = 1.034034135010395
$ws.Range("K12").Value = 1.031250253087389
$ws.Range("L12").Value = 1.03990970255773
$ws.Range("M12").Value = 1.048254121288018
$ws.Range("N12").Value = 1.01490796956697

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027419130380977
$ws.Range("D13").Value = 1.027627444782092
$ws.Range("E13").Value = 1.036382571445489
$ws.Range("F13").Value = 1.044761752456881
$ws.Range("I13").Value = 1.031003354329884
$ws.Range("J13").Value = 1.034109882441538
$ws.Range("K13").Value = 1.031261532100579
$ws.Range("L13").Value = 1.039983576638822
$ws.Range("M13").Value = 1.04833167048139
$ws.Range("N13").Value = 1.01493467238154

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027729839656846
$ws.Range("D14").Value = 1.027696536313468
$ws.Range("E14").Value = 1.036657146926726
$ws.Range("F14").Value = 1.045048094975385
$ws.Range("I14").Value = 1.031017899928756
$ws.Range("J14").Value = 1.034357975852198
$ws.Range("K14").Value = 1.031298452990928
$ws.Range("L14").Value = 1.040225534148502
$ws.Range("M14").Value = 1.048585660597373
$ws.Range("N14").Value = 1.015022112717112

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027921205470839
$ws.Range("D15").Value = 1.02773909712811
$ws.Range("E15").Value = 1.03682625929189
$ws.Range("F15").Value = 1.045224450073966
$ws.Range("I15").Value = 1.031026819027287
$ws.Range("J15").Value = 1.034510759126319
$ws.Range("K15").Value = 1.031321173887804
$ws.Range("L15").Value = 1.04037453861452
$ws.Range("M15").Value = 1.048742071732297
$ws.Range("N15").Value = 1.015075946794147

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029033944912591
$ws.Range("D16").Value = 1.027986684869852
$ws.Range("E16").Value = 1.037809620593889
$ws.Range("F16").Value = 1.04624985554159
$ws.Range("I16").Value = 1.031078076748075
$ws.Range("J16").Value = 1.035398885230365
$ws.Range("K16").Value = 1.031453002412254
$ws.Range("L16").Value = 1.041240696149573
$ws.Range("M16").Value = 1.049651232271605
$ws.Range("N16").Value = 1.015388665003065

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02973098954622
$ws.Range("D17").Value = 1.02814187063219
$ws.Range("E17").Value = 1.038425635902303
$ws.Range("F17").Value = 1.046892145899511
$ws.Range("I17").Value = 1.031109652355675
$ws.Range("J17").Value = 1.035954991921671
$ws.Range("K17").Value = 1.031535327386651
$ws.Range("L17").Value = 1.041783044570342
$ws.Range("M17").Value = 1.050220461233256
$ws.Range("N17").Value = 1.015584283209622

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030137228203359
$ws.Range("D18").Value = 1.028232344425954
$ws.Range("E18").Value = 1.038784655847776
$ws.Range("F18").Value = 1.04726645634995
$ws.Range("I18").Value = 1.031127861982604
$ws.Range("J18").Value = 1.03627900678409
$ws.Range("K18").Value = 1.031583213531832
$ws.Range("L18").Value = 1.04209904215984
$ws.Range("M18").Value = 1.050552103927997
$ws.Range("N18").Value = 1.015698190310197

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030275688688082
$ws.Range("D19").Value = 1.028263186248177
$ws.Range("E19").Value = 1.038907023382601
$ws.Range("F19").Value = 1.047394031574265
$ws.Range("I19").Value = 1.031134035737671
$ws.Range("J19").Value = 1.036389428066242
$ws.Range("K19").Value = 1.03159951896064
$ws.Range("L19").Value = 1.042206731040914
$ws.Range("M19").Value = 1.050665121647611
$ws.Range("N19").Value = 1.015736996976777

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029656238104468
$ws.Range("D20").Value = 1.028125225152969
$ws.Range("E20").Value = 1.038359573530083
$ws.Range("F20").Value = 1.046823268114908
$ws.Range("I20").Value = 1.031106286096164
$ws.Range("J20").Value = 1.035895363525731
$ws.Range("K20").Value = 1.031526508410133
$ws.Range("L20").Value = 1.041724891505313
$ws.Range("M20").Value = 1.05015942764409
$ws.Range("N20").Value = 1.015563315295731

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027638358639753
$ws.Range("D21").Value = 1.027676192402007
$ws.Range("E21").Value = 1.036576304368801
$ws.Range("F21").Value = 1.044963788845046
$ws.Range("I21").Value = 1.031013625555041
$ws.Range("J21").Value = 1.034284934231095
$ws.Range("K21").Value = 1.031287586391663
$ws.Range("L21").Value = 1.040154299036972
$ws.Range("M21").Value = 1.048510883645306
$ws.Range("N21").Value = 1.014996372225778

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026366489979174
$ws.Range("D22").Value = 1.027393484345619
$ws.Range("E22").Value = 1.035452368410141
$ws.Range("F22").Value = 1.043791616879112
$ws.Range("I22").Value = 1.030953492124369
$ws.Range("J22").Value = 1.033269118573119
$ws.Range("K22").Value = 1.031136175051535
$ws.Range("L22").Value = 1.039163604980436
$ws.Range("M22").Value = 1.047470872130625
$ws.Range("N22").Value = 1.01463813504079

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027041036019412
$ws.Range("D23").Value = 1.027543389189893
$ws.Range("E23").Value = 1.036048451147053
$ws.Range("F23").Value = 1.044413300826365
$ws.Range("I23").Value = 1.030985547694333
$ws.Range("J23").Value = 1.033807936813946
$ws.Range("K23").Value = 1.031216553865276
$ws.Range("L23").Value = 1.039689098454886
$ws.Range("M23").Value = 1.048022538809366
$ws.Range("N23").Value = 1.014828213436907

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029690016100355
$ws.Range("D24").Value = 1.028132746663079
$ws.Range("E24").Value = 1.038389425172652
$ws.Range("F24").Value = 1.04685439203021
$ws.Range("I24").Value = 1.031107807806412
$ws.Range("J24").Value = 1.035922308121047
$ws.Range("K24").Value = 1.031530493736001
$ws.Range("L24").Value = 1.041751169438359
$ws.Range("M24").Value = 1.050187007266663
$ws.Range("N24").Value = 1.015572790391942

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032750903629794
$ws.Range("D25").Value = 1.028814954861362
$ws.Range("E25").Value = 1.041094615144331
$ws.Range("F25").Value = 1.049674379042965
$ws.Range("I25").Value = 1.031241529135349
$ws.Range("J25").Value = 1.03836212856619
$ws.Range("K25").Value = 1.031889599187807
$ws.Range("L25").Value = 1.044130597398971
$ws.Range("M25").Value = 1.052683934115078
$ws.Range("N25").Value = 1.016429248683781
